$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on 2024-01-19.
# Column D holds "price" values formatted as plain text (e.g. thousand-dot
# separated numbers); force Text number format first so Excel does not
# reinterpret numeric-looking strings (e.g. "1.00", "311.49") as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.391.43"
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.467.30"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.49"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.51"
$ws.Range("E6").Value = "  -6.20%  "
$ws.Range("E7").Value = "  -3.43%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -4.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.49"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.01"
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.445.39"
$ws.Range("E15").Value = "  -4.42%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.98"
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("E17").Value = "  -3.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.350.27"
$ws.Range("E18").Value = "  -3.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.33"
$ws.Range("E19").Value = "  -5.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0924"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("E21").Value = "  -9.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.67"
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.75"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("E24").Value = "  -4.51%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -7.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.04"
$ws.Range("E27").Value = "  -6.09%  "
$ws.Range("E28").Value = "  -4.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.67"
$ws.Range("E29").Value = "  -5.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.46"
$ws.Range("E30").Value = "  -6.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.68"
$ws.Range("E31").Value = "  -5.08%  "
$ws.Range("E32").Value = "  -6.02%  "
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0749"
$ws.Range("E35").Value = "  -5.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.06"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.11"
$ws.Range("E37").Value = "  -6.60%  "
$ws.Range("E38").Value = "  -5.42%  "
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("E40").Value = "  -7.86%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("E43").Value = "  -9.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.989.07"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  -4.70%  "
$ws.Range("E46").Value = "  -9.20%  "
$ws.Range("E47").Value = "  -5.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.711.42"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "69.47"
$ws.Range("E49").Value = "  -4.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "96.42"
$ws.Range("E50").Value = "  -4.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.65"
$ws.Range("E51").Value = "  -6.66%  "
